{"js": "// Fix \"Tarefa 2\" intro sentence: \"...primeira equa\u00e7\u00e3o desta tarefa:\"\n// becomes \"...primeira equa\u00e7\u00e3o da tarefa 1:\"\nconst searchResults = context.document.body.search(\"desta tarefa:\", { matchCase: true });\nsearchResults.load(\"text\");\nawait context.sync();\n\nif (searchResults.items.length > 0) {\n  searchResults.items[0].insertText(\"da tarefa 1:\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Fix \"Tarefa 2\" intro sentence: \"...primeira equa\u00e7\u00e3o desta tarefa:\"\n# becomes \"...primeira equa\u00e7\u00e3o da tarefa 1:\"\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$found = $find.Execute(\"desta tarefa:\", $false, $false, $false, $false, $false, $true, 0, $false, \"da tarefa 1:\", 1)\n\nif (-not $found) {\n    Write-Output \"WARNING: target text 'desta tarefa:' not found\"\n}\n"}
